# Fix branch revenue export template: add a new "Tổng giảm giá đơn dưới 2.000đ"
# column to the revenue report header row (inserted just before the existing
# "Tổng xu" column), and move the active selection down near the bottom of
# the sheet (matching the author's post-edit cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (the "Tổng xu" column), shifting the
# remaining totals columns one place to the right and extending the header
# row's merged title cell (A1:K1 -> A1:L1) and sheet dimension accordingly.
$ws.Columns.Item(7).Insert()

# Populate the header for the newly inserted column.
$ws.Cells.Item(8, 7).Value = "Tổng giảm giá đơn dưới 2.000đ"

# Reflect the author's final cursor position after the edit.
$null = $ws.Range("G16").Select()
